$wb = $excel.ActiveWorkbook

# The fiscal year needs to be added to the "Global" sheet, next to the
# existing BrowserName / URL / dtFirstReqID columns, so it can be picked
# up (and seen) when costs are copied over from the other project.
$global = $wb.Worksheets.Item("Global")
$ppm = $wb.Worksheets.Item("PPMFinanceReview")

$global.Activate()

# D2 becomes the new "last column" of the row, so give it the boxed-edge
# border that used to belong to C2 (the previous last column).
$global.Range("C2").Copy()
$global.Range("D2").PasteSpecial(-4122) # xlPasteFormats

# C2 is no longer the last column -- it now looks like the other inner
# data cells (A2/B2), i.e. no right-hand border.
$global.Range("A2").Copy()
$global.Range("C2").PasteSpecial(-4122) # xlPasteFormats

$global.Range("D1").Value = "FiscalYear"
$global.Range("D2").Value = 2020

# Give the new column its own (slightly narrower) width instead of sharing
# the sheet's blanket default.
$global.Columns.Item(4).ColumnWidth = 8.58

$global.Range("D2").Select()

# Restore the originally-active sheet/tab (PPMFinanceReview).
$ppm.Activate()
